# Apply updated object-detection / bounding-box results to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The bounding_box_confidence column (K) holds numeric-looking values that
# are actually stored as text in the workbook; force text formatting so the
# assignments below don't get auto-coerced to numbers.
$ws.Range("K2:K16").NumberFormat = "@"

# Row 2
$ws.Range("B2").Value = "poster"
$ws.Range("J2").Value = "[    0.50035     0.50007     0.98906     0.99065]"
$ws.Range("K2").Value = "0.3906828"

# Row 3
$ws.Range("B3").Value = "snow"
$ws.Range("J3").Value = "[     0.5001     0.78907           1      0.4214]"
$ws.Range("K3").Value = "0.4067031"

# Row 4
$ws.Range("B4").Value = "glacier"
$ws.Range("J4").Value = "[    0.50144     0.47236     0.23467     0.31704]"
$ws.Range("K4").Value = "0.38370985"

# Row 5
$ws.Range("B5").Value = "person"
$ws.Range("J5").Value = "[    0.55918     0.55452     0.47884     0.57723]"
$ws.Range("K5").Value = "0.47672573"

# Row 6
$ws.Range("B6").Value = "tree,podium,plant"
$ws.Range("J6").Value = "[    0.48386     0.45351     0.31128      0.3735],[    0.72244     0.61256     0.16141     0.14683],[    0.48364     0.45348     0.31012     0.37298]"
$ws.Range("K6").Value = "0.6849565, 0.5135502, 0.38693547"

# Row 7
$ws.Range("B7").Value = "podium,person,flag"
$ws.Range("J7").Value = "[    0.18573     0.60911       0.121     0.24808],[    0.64521     0.47421     0.58957     0.71254],[    0.18102     0.42178    0.072017     0.15551]"
$ws.Range("K7").Value = "0.50541675, 0.4737043, 0.42418364"

# Row 8
$ws.Range("B8").Value = "glacier,sign,camera,map,mountain,poster"
$ws.Range("J8").Value = "[    0.25408     0.51953     0.50421     0.94968],[    0.25361     0.52039     0.50393     0.94866],[    0.84931     0.57744    0.038343    0.066539],[    0.25352     0.52051     0.50352     0.94984],[    0.74161      0.2068     0.47352     0.10733],[    0.25401     0.52041     0.50412      0.9485]"
$ws.Range("K8").Value = "0.39810526, 0.5070563, 0.41183785, 0.43476003, 0.4302195, 0.45908466"

# Row 9
$ws.Range("B9").Value = "person,book,poster,podium"
$ws.Range("J9").Value = "[    0.48141     0.56795     0.58878       0.858],[    0.12911     0.79567     0.25538      0.2004],[    0.42391     0.18381     0.55161     0.36553],[    0.15456       0.798     0.30594     0.39658]"
$ws.Range("K9").Value = "0.4189234, 0.39079067, 0.41635087, 0.37884402"

# Row 10 (object detection found nothing here any more - clear the cells,
# keeping them present-but-blank like the rest of the sheet's empty cells)
$ws.Range("B10").NumberFormat = "@"
$ws.Range("J10").NumberFormat = "@"
$ws.Range("B10").Value = ""
$ws.Range("J10").Value = ""
$ws.Range("K10").Value = ""

# Row 11
$ws.Range("B11").Value = "person"
$ws.Range("J11").Value = "[    0.21822     0.69391     0.43636     0.61248]"
$ws.Range("K11").Value = "0.6345378"

# Row 12
$ws.Range("B12").Value = "map"
$ws.Range("J12").Value = "[    0.81969     0.34883      0.2358     0.60134]"
$ws.Range("K12").Value = "0.4421639"

# Row 13
$ws.Range("B13").Value = "person,tree"
$ws.Range("J13").Value = "[    0.21239     0.66656     0.29598     0.31358],[    0.50171     0.41076     0.98706     0.81836]"
$ws.Range("K13").Value = "0.51314807, 0.39316073"

# Row 14
$ws.Range("B14").Value = "person,flower"
$ws.Range("J14").Value = "[    0.29658      0.5592     0.58964     0.87699],[     0.9504     0.49174     0.09288    0.093094]"
$ws.Range("K14").Value = "0.43543482, 0.37039283"

# Row 15
$ws.Range("B15").Value = "sign,person"
$ws.Range("J15").Value = "[     0.6372     0.03649     0.12775    0.051663],[    0.16469     0.49228     0.30822     0.43286]"
$ws.Range("K15").Value = "0.39736262, 0.42133933"

# Row 16
$ws.Range("B16").Value = "car"
$ws.Range("J16").Value = "[    0.32883     0.73569     0.62881     0.24144]"
$ws.Range("K16").Value = "0.40813386"
